# "Generate Report for Archive"
#
# The localization status "Ready for handoff" has moved on to
# "In Translation" for every locale tracked in this report. Refresh the
# status text wherever it appears (the Overview sheet shows one status
# column per locale, and each per-locale sheet - zh-cn, de-de - has its own
# Status column), then resize those status columns to fit the new, shorter
# text, the way a regenerated report naturally would.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Replace every occurrence of the old status text with the new one, on
# every sheet, instead of hard-coding cell addresses.
foreach ($ws in $wb.Worksheets) {
    $replaced = $ws.Cells.Replace($oldStatus, $newStatus)
}

# Width (in "characters") that the status columns should shrink to, now
# that the longest value they hold is the shorter "In Translation" instead
# of "Ready for handoff".
$statusColumnWidth = 12.5

# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $statusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColumnWidth

# Per-locale detail sheets: column C holds the Status value.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColumnWidth
